$d = $word.ActiveDocument

# The table cell (right column, under the Facebook hyperlink) currently has:
#   - an empty paragraph
#   - a paragraph with the sentence "Esta página é importante, pois permite
#     a divulgação e marketing dos serviços prestados." which also carries
#     a "_GoBack" bookmark right before the final period.
# The edit removes the empty paragraph and strips the sentence's wording,
# leaving a single empty paragraph that still holds the "_GoBack" bookmark.

$searchText = "Esta página é importante*prestados."

function Find-SentenceRange($doc) {
    $r = $doc.Content.Duplicate
    $ok = $r.Find.Execute($searchText, $true, $false, $true, $false, $false, `
        $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not locate the target sentence."
    }
    return $r
}

# --- Step 1: remove the empty paragraph right before the sentence paragraph
$sentenceRange = Find-SentenceRange $d
$sentencePara = $sentenceRange.Paragraphs(1)
$sentenceParaStart = $sentencePara.Range.Start
if ($sentenceParaStart -gt 0) {
    $prevRange = $d.Range($sentenceParaStart - 1, $sentenceParaStart - 1)
    $prevPara = $prevRange.Paragraphs(1)
    if ($prevPara.Range.Text -eq "`r") {
        $prevPara.Range.Delete()
    }
}

# --- Step 2: clear the sentence's wording while preserving the bookmark
# that sits between "...prestados" and the final ".". Deleting the text
# before the bookmark and the text after it as two separate operations
# (instead of one delete spanning the whole sentence) keeps the bookmark
# alive, since neither delete range fully encloses the bookmark itself.
$sentenceRange = Find-SentenceRange $d
$sentencePara = $sentenceRange.Paragraphs(1)
$paraStart = $sentencePara.Range.Start
$paraEnd = $sentencePara.Range.End

# Delete the trailing "." (the character right before the paragraph mark).
$d.Range($paraEnd - 2, $paraEnd - 1).Delete()

# Delete the leading wording, up to (not including) the bookmark/period spot.
$d.Range($paraStart, $paraEnd - 2).Delete()
